$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to seed the shared-string pool with a clean text
# representation of numeric-looking values before writing them into the
# target cell with a leading apostrophe (keeps the cell text-typed without
# introducing a new persistent number format on the scratch cell itself).
$scratch = $ws.Range("Z100")

function Set-TextValue($cell, [string]$text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $cell.Value = "'" + $text
    $scratch.Clear()
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '35.535.48'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.913.20'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  -0.18%  '
Set-TextValue $ws.Range("D5") '0.704'
$ws.Range("E5").Value = '  +6.64%  '
Set-TextValue $ws.Range("D6") '247.51'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("E7").Value = '  -0.15%  '
Set-TextValue $ws.Range("D8") '40.73'
$ws.Range("E8").Value = '  -3.17%  '
Set-TextValue $ws.Range("D9") '0.358'
$ws.Range("E9").Value = '  +3.87%  '
Set-TextValue $ws.Range("D10") '52.76'
$ws.Range("E10").Value = '  +6.54%  '
Set-TextValue $ws.Range("D11") '0.0732'
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("E13").Value = '  +0.17%  '
Set-TextValue $ws.Range("D14") '12.74'
$ws.Range("E14").Value = '  +3.13%  '
Set-TextValue $ws.Range("D15") '0.718'
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").Value = '1.914.93'
$ws.Range("E16").Value = '  +0.24%  '
Set-TextValue $ws.Range("D17") '4.92'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = '35.519.55'
$ws.Range("E18").Value = '  -0.07%  '
Set-TextValue $ws.Range("D19") '73.34'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("E20").Value = '  -0.58%  '
Set-TextValue $ws.Range("D21") '13.18'
$ws.Range("E21").Value = '  +3.97%  '
Set-TextValue $ws.Range("D22") '242.70'
$ws.Range("E22").Value = '  -0.79%  '
Set-TextValue $ws.Range("D23") '5.08'
$ws.Range("E23").Value = '  +4.59%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +1.11%  '
Set-TextValue $ws.Range("D26") '2.30'
$ws.Range("E26").Value = '  +3.47%  '
Set-TextValue $ws.Range("D27") '168.94'
$ws.Range("E27").Value = '  -1.57%  '
Set-TextValue $ws.Range("D28") '8.65'
$ws.Range("E28").Value = '  +2.04%  '
Set-TextValue $ws.Range("D29") '18.81'
$ws.Range("E29").Value = '  +1.79%  '
$ws.Range("E30").Value = '  +2.64%  '
$ws.Range("D31").Value = '4.142.65'
$ws.Range("E31").Value = '  +19.86%  '
$ws.Range("E32").Value = '  +1.57%  '
Set-TextValue $ws.Range("D33") '0.0577'
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("E34").Value = '  +10.23%  '
Set-TextValue $ws.Range("D35") '4.22'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("E36").Value = '  -0.16%  '
Set-TextValue $ws.Range("D37") '0.909'
$ws.Range("E37").Value = '  -5.72%  '
$ws.Range("E38").Value = '  +11.39%  '
$ws.Range("E39").Value = '  +0.79%  '
Set-TextValue $ws.Range("D40") '17.43'
$ws.Range("E40").Value = '  +11.47%  '
Set-TextValue $ws.Range("D41") '98.61'
$ws.Range("E41").Value = '  +6.90%  '
$ws.Range("E42").Value = '  +2.46%  '
Set-TextValue $ws.Range("D43") '0.0211'
$ws.Range("E43").Value = '  +2.38%  '
Set-TextValue $ws.Range("D44") '0.0650'
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").Value = '1.352.44'
$ws.Range("E45").Value = '  +0.36%  '
Set-TextValue $ws.Range("D46") '2.46'
$ws.Range("E46").Value = '  +2.54%  '
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("E48").Value = '  +0.24%  '
Set-TextValue $ws.Range("D49") '45.70'
$ws.Range("E49").Value = '  -3.81%  '
Set-TextValue $ws.Range("D50") '12.22'
$ws.Range("E50").Value = '  -3.55%  '
Set-TextValue $ws.Range("D51") '6.56'
$ws.Range("E51").Value = '  -0.21%  '
